$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H108").Value = 30000
$ws.Range("J108").Value = 30000
$ws.Range("L108").Value = 30000
$ws.Range("N108").Value = -37680
$ws.Range("H132").Value = 1565.4791
$ws.Range("I132").Value = 1452.5217
$ws.Range("J132").Value = 4163.5
$ws.Range("K132").Value = 4357.5651
$ws.Range("L132").Value = 12490.5
$ws.Range("M132").Value = -1827.5651
$ws.Range("N132").Value = -17550.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2946.9167
$ws.Range("I2").Value = 1213.8334
$ws.Range("J2").Value = 4680
$ws.Range("K2").Value = 1213.8334
$ws.Range("L2").Value = 4680
$ws.Range("M2").Value = -1100.8334
$ws.Range("N2").Value = -4906
$ws.Range("H74").Value = 1202
$ws.Range("I74").Value = 712.0952
$ws.Range("K74").Value = 712.0952
$ws.Range("M74").Value = 161.9048
$ws.Range("H77").Value = 1202
$ws.Range("I77").Value = 712.0952
$ws.Range("K77").Value = 3560.476
$ws.Range("M77").Value = 807.5240000000003
$ws.Range("H97").Value = 1410.9722
$ws.Range("I97").Value = 1197.6522
$ws.Range("J97").Value = 1788.3846
$ws.Range("K97").Value = 1197.6522
$ws.Range("L97").Value = 1788.3846
$ws.Range("M97").Value = -701.6522
$ws.Range("N97").Value = -2780.3846
$ws.Range("H107").Value = 0
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()
$ws.Range("H108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680
$ws.Range("H110").Value = 2064.8667
$ws.Range("I110").Value = 1707.1
$ws.Range("J110").Value = 2780.4
$ws.Range("K110").Value = 1707.1
$ws.Range("L110").Value = 2780.4
$ws.Range("M110").Value = 337.9000000000001
$ws.Range("N110").Value = -6870.4
$ws.Range("H111").Value = 35000
$ws.Range("J111").Value = 35000
$ws.Range("L111").Value = 35000
$ws.Range("N111").Value = -43180
$ws.Range("H112").Value = 10901.75
$ws.Range("J112").Value = 10901.75
$ws.Range("L112").Value = 10901.75
$ws.Range("N112").Value = -13855.75
$ws.Range("H113").Value = 50000
$ws.Range("J113").Value = 50000
$ws.Range("L113").Value = 50000
$ws.Range("N113").Value = -58678
$ws.Range("H114").Value = 11066.333
$ws.Range("J114").Value = 11066.333
$ws.Range("L114").Value = 11066.333
$ws.Range("N114").Value = -19744.333
$ws.Range("H116").Value = 2946.9167
$ws.Range("I116").Value = 1213.8334
$ws.Range("J116").Value = 4680
$ws.Range("K116").Value = 1213.8334
$ws.Range("L116").Value = 4680
$ws.Range("M116").Value = 1080.1666
$ws.Range("N116").Value = -9268
$ws.Range("H117").Value = 27600
$ws.Range("J117").Value = 27600
$ws.Range("L117").Value = 27600
$ws.Range("N117").Value = -36778
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()
$ws.Range("H119").Value = 20060.777
$ws.Range("J119").Value = 20060.777
$ws.Range("L119").Value = 20060.777
$ws.Range("N119").Value = -29736.777
$ws.Range("H120").Value = 27400
$ws.Range("J120").Value = 27400
$ws.Range("L120").Value = 27400
$ws.Range("N120").Value = -37076
$ws.Range("H121").Value = 27900
$ws.Range("J121").Value = 27900
$ws.Range("L121").Value = 27900
$ws.Range("N121").Value = -31394
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2946.9167
$ws.Range("I3").Value = 1213.8334
$ws.Range("J3").Value = 4680
$ws.Range("K3").Value = 1213.8334
$ws.Range("L3").Value = 4680
$ws.Range("M3").Value = -1099.8334
$ws.Range("N3").Value = -4908
$ws.Range("H20").Value = 2464.2
$ws.Range("I20").Value = 1832.0454
$ws.Range("J20").Value = 3236.8333
$ws.Range("K20").Value = 1832.0454
$ws.Range("L20").Value = 3236.8333
$ws.Range("M20").Value = -1585.0454
$ws.Range("N20").Value = -3730.8333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3180.2
$ws.Range("I16").Value = 2975.25
$ws.Range("J16").Value = 4000
$ws.Range("K16").Value = 2975.25
$ws.Range("L16").Value = 4000
$ws.Range("M16").Value = -2688.25
$ws.Range("N16").Value = -4574
$ws.Range("H113").Value = 3180.2
$ws.Range("I113").Value = 2975.25
$ws.Range("J113").Value = 4000
$ws.Range("K113").Value = 2975.25
$ws.Range("L113").Value = 4000
$ws.Range("M113").Value = -805.25
$ws.Range("N113").Value = -8340
$ws.Range("H122").Value = 1870.3125
$ws.Range("I122").Value = 1703.4783
$ws.Range("J122").Value = 2296.6667
$ws.Range("K122").Value = 5110.4349
$ws.Range("L122").Value = 6890.000100000001
$ws.Range("M122").Value = -2660.4349
$ws.Range("N122").Value = -11790.0001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 400
$ws.Range("I97").Value = 200
$ws.Range("J97").Value = 500
$ws.Range("K97").Value = 600
$ws.Range("L97").Value = 1500
$ws.Range("M97").Value = -104
$ws.Range("N97").Value = -2492
$ws.Range("H107").Value = 392.63635
$ws.Range("I107").Value = 200
$ws.Range("J107").Value = 411.9
$ws.Range("K107").Value = 600
$ws.Range("L107").Value = 1235.7
$ws.Range("M107").Value = 1320
$ws.Range("N107").Value = -5075.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H109").Value = 17213.75
$ws.Range("J109").Value = 17213.75
$ws.Range("L109").Value = 17213.75
$ws.Range("N109").Value = -19293.75
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H113").Value = 76925064
$ws.Range("I113").Value = 1992.2
$ws.Range("J113").Value = 125001976
$ws.Range("K113").Value = 1992.2
$ws.Range("L113").Value = 125001976
$ws.Range("M113").Value = 177.8
$ws.Range("N113").Value = -125006316
$ws.Range("H114").Value = 48000
$ws.Range("J114").Value = 48000
$ws.Range("L114").Value = 48000
$ws.Range("N114").Value = -56678
$ws.Range("H116").Value = 40000
$ws.Range("J116").Value = 40000
$ws.Range("L116").Value = 40000
$ws.Range("N116").Value = -49178
$ws.Range("H117").Value = 30000
$ws.Range("J117").Value = 30000
$ws.Range("L117").Value = 30000
$ws.Range("N117").Value = -36884
$ws.Range("H118").Value = 44990
$ws.Range("J118").Value = 44990
$ws.Range("L118").Value = 44990
$ws.Range("N118").Value = -48304
$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()
$ws.Range("H121").Value = 30000
$ws.Range("J121").Value = 30000
$ws.Range("L121").Value = 30000
$ws.Range("N121").Value = -33494
$ws.Range("H132").Value = 4405.62
$ws.Range("I132").Value = 4532.225
$ws.Range("K132").Value = 13596.675
$ws.Range("M132").Value = -11066.675

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2594.3333
$ws.Range("I61").Value = 1799.9
$ws.Range("J61").Value = 3587.375
$ws.Range("K61").Value = 1799.9
$ws.Range("L61").Value = 3587.375
$ws.Range("M61").Value = -1597.9
$ws.Range("N61").Value = -3991.375
$ws.Range("H108").Value = 0
$ws.Range("I108").Value = 0
$ws.Range("K108").Value = 0
$ws.Range("M108").ClearContents()
$ws.Range("H110").Value = 20240.666
$ws.Range("J110").Value = 20240.666
$ws.Range("L110").Value = 20240.666
$ws.Range("N110").Value = -28420.666
$ws.Range("H112").Value = 30000
$ws.Range("J112").Value = 30000
$ws.Range("L112").Value = 30000
$ws.Range("N112").Value = -32954
$ws.Range("H113").Value = 2594.3333
$ws.Range("I113").Value = 1799.9
$ws.Range("J113").Value = 3587.375
$ws.Range("K113").Value = 1799.9
$ws.Range("L113").Value = 3587.375
$ws.Range("M113").Value = 370.0999999999999
$ws.Range("N113").Value = -7927.375
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("L115").Value = 0
$ws.Range("N115").ClearContents()
$ws.Range("H116").Value = 20600
$ws.Range("J116").Value = 20600
$ws.Range("L116").Value = 20600
$ws.Range("N116").Value = -29778
$ws.Range("H117").Value = 25000
$ws.Range("J117").Value = 25000
$ws.Range("L117").Value = 25000
$ws.Range("N117").Value = -34178
$ws.Range("H119").Value = 29500
$ws.Range("J119").Value = 29500
$ws.Range("L119").Value = 29500
$ws.Range("N119").Value = -39176
$ws.Range("H120").Value = 32000
$ws.Range("J120").Value = 32000
$ws.Range("L120").Value = 32000
$ws.Range("N120").Value = -41676
$ws.Range("H121").Value = 20616.666
$ws.Range("J121").Value = 20616.666
$ws.Range("L121").Value = 20616.666
$ws.Range("N121").Value = -24110.666
$ws.Range("H122").Value = 2105.3438
$ws.Range("I122").Value = 1927.76
$ws.Range("J122").Value = 2739.5715
$ws.Range("K122").Value = 5783.28
$ws.Range("L122").Value = 8218.7145
$ws.Range("M122").Value = -3333.28
$ws.Range("N122").Value = -13118.7145
